# Update cryptocurrency price/volume data on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

# Row 2 - Bitcoin
Set-Text "D2" "26.457.47"
Set-Text "E2" "  -0.74%  "

# Row 3 - Ethereum
Set-Text "D3" "1.622.79"
Set-Text "E3" "  -0.06%  "

# Row 4 - TetherUSD
Set-Text "E4" "  +0.38%  "

# Row 5 - BNB
Set-Text "D5" "213.38"
Set-Text "E5" "  -0.70%  "

# Row 6 - XRP
Set-Text "D6" "0.501"
Set-Text "E6" "  -0.90%  "

# Row 7 - USDC
Set-Text "E7" "  +0.37%  "

# Row 8 - Cardano
Set-Text "D8" "0.246"
Set-Text "E8" "  -0.23%  "

# Row 9 - Dogecoin
Set-Text "D9" "0.0609"
Set-Text "E9" "  -0.43%  "

# Row 10 - Solana
Set-Text "D10" "19.15"
Set-Text "E10" "  -1.25%  "

# Row 11 - TRON
Set-Text "D11" "0.0855"
Set-Text "E11" "  -0.17%  "

# Row 12 - Wrapped liquid staked Ether 2.0
Set-Text "D12" "1.854.28"
Set-Text "E12" "  +0.20%  "

# Row 13 - Wrapped Ether
Set-Text "D13" "1.629.87"
Set-Text "E13" "  -0.17%  "

# Row 14 - Polkadot
Set-Text "E14" "  -0.37%  "

# Row 15 - Polygon
Set-Text "D15" "0.511"
Set-Text "E15" "  -0.64%  "

# Row 16 - Litecoin
Set-Text "E16" "  -1.64%  "

# Row 17 - Bitcoin Cash
Set-Text "D17" "234.32"
Set-Text "E17" "  +1.08%  "

# Row 18 - Wrapped BTC
Set-Text "D18" "26.484.86"
Set-Text "E18" "  -0.52%  "

# Row 19 - Chainlink
Set-Text "D19" "7.76"
Set-Text "E19" "  +1.70%  "

# Row 20 - Shiba Inu
Set-Text "D20" "0.0₃0725"
Set-Text "E20" "  -0.57%  "

# Row 21 - Dai
Set-Text "E21" "  +0.29%  "

# Row 22 - Uniswap
Set-Text "E22" "  -2.05%  "

# Row 23 - Toncoin
Set-Text "E23" "  -1.78%  "

# Row 24 - Avalanche
Set-Text "D24" "9.13"
Set-Text "E24" "  -0.07%  "

# Row 25 - Monero
Set-Text "D25" "146.92"
Set-Text "E25" "  +0.97%  "

# Row 26 - BinanceUSD
Set-Text "E26" "  +0.34%  "

# Row 27 - Cosmos
Set-Text "D27" "7.06"
Set-Text "E27" "  -0.08%  "

# Row 28 - Stellar
Set-Text "E28" "  -0.69%  "

# Row 29 - Ethereum Classic
Set-Text "D29" "15.61"
Set-Text "E29" "  -0.30%  "

# Row 30 - Hedera
Set-Text "E30" "  -0.57%  "

# Row 31 - PancakeSwap
Set-Text "E31" "  -0.38%  "

# Row 32 - Maker
Set-Text "D32" "1.515.64"
Set-Text "E32" "  +4.75%  "

# Row 33 - Filecoin
Set-Text "D33" "3.25"
Set-Text "E33" "  -0.07%  "

# Row 34 - Internet Computer (DFINITY)
Set-Text "E34" "  -1.05%  "

# Row 35 - Lido DAO Token
Set-Text "E35" "  +2.14%  "

# Row 36 - Huobi Token
Set-Text "E36" "  +0.21%  "

# Row 37 - Immutable X
Set-Text "D37" "0.566"
Set-Text "E37" "  +0.99%  "

# Row 38 - VeChain
Set-Text "E38" "  -0.79%  "

# Row 39 - ARBITRUM
Set-Text "D39" "0.833"
Set-Text "E39" "  -0.75%  "

# Row 40 - FraxShare
Set-Text "E40" "  -0.27%  "

# Row 41 - PaxDollar
Set-Text "E41" "  +0.29%  "

# Row 42 - MXToken
Set-Text "E42" "  -0.13%  "

# Row 43 - RocketPoolETH
Set-Text "D43" "1.764.62"
Set-Text "E43" "  +0.18%  "

# Row 44 - Aave
Set-Text "D44" "62.77"
Set-Text "E44" "  +0.88%  "

# Row 45 - TrustWalletToken
Set-Text "D45" "0.761"
Set-Text "E45" "  -0.30%  "

# Row 46 - WEMIXToken
Set-Text "D46" "0.911"
Set-Text "E46" "  -3.45%  "

# Row 47 - Quant
Set-Text "D47" "89.70"
Set-Text "E47" "  +1.26%  "

# Row 48 - RenderToken
Set-Text "E48" "  -0.21%  "

# Row 49 - was BabyDogeCoin, now Cronos
Set-Text "B49" "Cronos"
Set-Text "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-Text "D49" "0.0501"
Set-Text "E49" "  -0.50%  "

# Row 50 - was Cronos, now EnergySwap
Set-Text "B50" "EnergySwap"
Set-Text "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Text "D50" "7.54"
Set-Text "E50" "  +0.65%  "

# Row 51 - was EnergySwap, now Algorand
Set-Text "B51" "Algorand"
Set-Text "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-Text "D51" "0.0962"
Set-Text "E51" "  -0.65%  "
